$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new data row: username "dpaul" / password "password"
$ws.Range("A3").Value = "dpaul"
$ws.Range("B3").Value = "password"

# Update the saved selection to A4 (matches target diff)
$ws.Range("A4").Select()
